$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Organograma")
$ws.Activate()

# Insert a new blank row at row 18 (shifts rows 18+ down by one)
$ws.Rows.Item(18).Insert()

# Fill in the new last row (row 36) with the new team member
$ws.Cells.Item(36, 2).Value = "Key User N1"
$ws.Cells.Item(36, 3).Value = "Luciana de Carvalho Cavalcante"
$ws.Cells.Item(36, 4).Value = "Coordenadora Juridico  "

Write-Output "done"
